# Copiar fecha disponible y agregarla al grupo que necesita; termine el
# flujo cuando se agrega un nuevo grupo al excel y se crea la carpeta
# para continuar el flujo.
$wb = $excel.ActiveWorkbook

$wsGrupos = $wb.Worksheets.Item("grupos")
$wsFechas = $wb.Worksheets.Item("fechas")

# "Grupo2" (fila 3 de "grupos") necesita una fecha disponible: la
# siguiente fecha libre en "fechas" es la de la fila 3 (A3), que todavia
# no tiene marca de "Uso". La copiamos hacia la hoja de grupos junto con
# la nueva cantidad.
$wsGrupos.Range("B3").Value = 5

$wsFechas.Range("A3").Copy() | Out-Null
$wsGrupos.Range("C3").PasteSpecial(-4163) | Out-Null  # xlPasteValues

# Marcamos esa fecha como usada ("x") en la hoja "fechas" y dejamos la
# seleccion sobre la celda recien usada.
$wsFechas.Range("B3").Value = "x"
$wsFechas.Range("A3").Select() | Out-Null

# El flujo termina seleccionando la siguiente fila vacia de "grupos"
# para poder agregar el nuevo grupo.
$wsGrupos.Rows("4:4").Select() | Out-Null

$wb.Save() | Out-Null
